$wb = $excel.ActiveWorkbook

# --- Update the existing "Contacts" sheet (sheet1) ---
$contacts = $wb.Worksheets.Item("Contacts")

# Header row changes: LifeCycleStage -> StageName, add LeadStatus
$contacts.Range("D1").Value = "StageName"
$contacts.Range("E1").Value = "LeadStatus"

# Data row changes
$contacts.Range("A2").Value = "ab2@abc.com"
$contacts.Range("B2").Value = "Razz"
$contacts.Range("C2").Value = "Ram"
$contacts.Range("E2").Value = "New"

# --- Add the new "Deals" sheet (sheet2) after Contacts ---
$deals = $wb.Worksheets.Add($null, $contacts)
$deals.Name = "Deals"

$deals.Range("A1:C1").ColumnWidth = 15.5

$deals.Range("A1").Value = "DealName"
$deals.Range("B1").Value = "DealMonth"
$deals.Range("C1").Value = "DealDay"
$deals.Range("D1").Value = "DealAmount"

$deals.Range("A2").Value = "Test"
$deals.Range("B2").NumberFormat = "mmmm\ yyyy"
$deals.Range("B2").Value = (Get-Date -Year 2020 -Month 6 -Day 1 -Hour 0 -Minute 0 -Second 0)
$deals.Range("C2").Value = 25
$deals.Range("D2").Value = 20

$null = $deals.Range("D14").Select()
